$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R2").Value = $null

$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 359547
$ws.Range("R5").Value = 5.3
$ws.Range("R6").Value = 99075
$ws.Range("R7").Value = 1.5
